# Update E_Magnitude values for specific rows (reduce by 0.20)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(6, 9, 11, 14, 16, 19, 21, 24, 26, 29, 31, 34, 36, 39, 41, 44, 46, 49, 51, 54, 56, 61)

foreach ($r in $rows) {
    $cell = $ws.Range("E$r")
    $current = $cell.Value2
    $cell.Value2 = [math]::Round($current - 0.20, 2)
}
